$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Case")

# --- New "Level Classification" / "A" row (E12/F12) ---
# E12 needs the bordered/shaded label look used by the other label cells in
# this column (e.g. E10/E11), so copy that formatting over before writing
# the text; F12 keeps its existing (plain) formatting.
$ws.Range("E10").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("E12").Value = "Level Classification"
$ws.Range("F12").Value = "A"

# --- Step 3 text updates (E16/F16) ---
$ws.Range("E16").Value = "Press Tab (or Shift+Tab to traverse reverse)from keyboard and navigate inside the website; Use down/up/left/right arrow keys to navigate to the images available in the page" + [char]10 + "Note: Avoid using mouse from this step; Keyboard usage is recommended"

$ws.Range("F16").Value = "Screen reader should read in proper sequence such as top left to bottom right when you navigate." + [char]10 + "Note: If you keep pressing TAB to move to different elements, some elements let you to navigate down and when pressing next TAB result in top most section of the page (instead going down to next item),it can be considered as confused structure or defect in order to get it fixed for a clear flow of sequence. " + [char]10 + "Verification Step: Follow the same sequence in reverse by using SHIFT+TAB to make sure that the reverse sequence is going back to the elements appropriate."

# Row grew taller to fit the longer verification text.
$ws.Rows.Item(16).RowHeight = 169.5

# --- View / selection state ---
$ws.Range("F18").Select()
